$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ticket category" (Categoria de Ingresso) requirement row first.
# This is originally row 5 ("RF-04 | Exibir Categoria de Ingresso").
$ws.Rows.Item(5).Delete()

# After the above deletion, the remaining "Categoria de Ingresso" CRUD rows
# (originally RF-14/RF-15/RF-16) have shifted up to rows 14-16.
# Select that block and delete it, matching the author's final selection state.
$ws.Range("A14:XFD16").Select()
$ws.Range("A14:A16").EntireRow.Delete()
